$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - F column "想去人数" (want-to-go count) updates
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F7").Value  = 4459
$wsExhibit.Range("F8").Value  = 2647
$wsExhibit.Range("F10").Value = 2589
$wsExhibit.Range("F11").Value = 1951
$wsExhibit.Range("F13").Value = 1678
$wsExhibit.Range("F14").Value = 688
$wsExhibit.Range("F15").Value = 140
$wsExhibit.Range("F16").Value = 159
$wsExhibit.Range("F18").Value = 31
$wsExhibit.Range("F20").Value = 77
$wsExhibit.Range("F21").Value = 42
$wsExhibit.Range("F23").Value = 32
$wsExhibit.Range("F25").Value = 580
$wsExhibit.Range("F29").Value = 440
$wsExhibit.Range("F31").Value = 1179
$wsExhibit.Range("F32").Value = 200
$wsExhibit.Range("F33").Value = 26
$wsExhibit.Range("F34").Value = 1256
$wsExhibit.Range("F35").Value = 2131
$wsExhibit.Range("F36").Value = 306
$wsExhibit.Range("F37").Value = 17
$wsExhibit.Range("F39").Value = 93
$wsExhibit.Range("F41").Value = 83
$wsExhibit.Range("F42").Value = 35
$wsExhibit.Range("F43").Value = 695
$wsExhibit.Range("F44").Value = 1368
$wsExhibit.Range("F45").Value = 132
$wsExhibit.Range("F48").Value = 57
$wsExhibit.Range("F49").Value = 81

# Sheet "演出" (Performance) - F column update
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F12").Value = 12

# Sheet "全部类型" (All Types) - F column "想去人数" updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 4459
$wsAll.Range("F6").Value  = 2647
$wsAll.Range("F7").Value  = 2589
$wsAll.Range("F8").Value  = 1678
$wsAll.Range("F11").Value = 688
$wsAll.Range("F12").Value = 140
$wsAll.Range("F13").Value = 159
$wsAll.Range("F15").Value = 31
$wsAll.Range("F17").Value = 77
$wsAll.Range("F19").Value = 32
$wsAll.Range("F20").Value = 580
$wsAll.Range("F27").Value = 440
$wsAll.Range("F29").Value = 1179
$wsAll.Range("F30").Value = 200
$wsAll.Range("F33").Value = 2131
$wsAll.Range("F34").Value = 306
$wsAll.Range("F37").Value = 12
$wsAll.Range("F39").Value = 93
$wsAll.Range("F41").Value = 83
$wsAll.Range("F42").Value = 35
$wsAll.Range("F43").Value = 695
$wsAll.Range("F44").Value = 1368
$wsAll.Range("F46").Value = 132
$wsAll.Range("F48").Value = 81
